# Apply the "routes" schema sheet update:
# - insert a new row for the "person" column definition
# - change approverLayer's Type from tinyint(1) to int(5)
# - rename "title" to "approverTitle"
# - refresh column widths / selection to match the author's final view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 4 (shifts existing rows 4-19 down to 5-20)
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the "person" column (only ColumnID + Type are known so far)
$ws.Cells.Item(4, 1).Value = "person"
$ws.Cells.Item(4, 3).Value = "varchar(255)"

# approverLayer (now row 9) Type corrected from tinyint(1) to int(5)
$ws.Cells.Item(9, 3).Value = "int(5)"

# "title" column (now row 11) renamed to "approverTitle"
$ws.Cells.Item(11, 1).Value = "approverTitle"

# Column widths adjusted by the author
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 23.833333333333332
$ws.Columns.Item(3).ColumnWidth = 14.333333333333334

# Selection moved to A12
$ws.Range("A12").Select() | Out-Null
